# Generate Report for Handoff
# The f5ba21a8-16a3-4f7b-a38c-01b440ee532d.md file moved from "In Translation"
# to "Ready for handoff": update the Overview, zh-cn and de-de sheets for the
# row that corresponds to that file (row 3 on every sheet).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-09-05 08:18:23"

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "mt"
$zhcn.Range("G3").Value = "f5ba21a8-16a3-4f7b-a38c-01b440ee532d.3458621fa1e454d0c63bdc08382b6425b2e94f4c.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-09-05 08:18:19"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "mt"
$dede.Range("G3").Value = "f5ba21a8-16a3-4f7b-a38c-01b440ee532d.3458621fa1e454d0c63bdc08382b6425b2e94f4c.de-de.xlf"
$dede.Range("H3").Value = "2016-09-05 08:18:23"

# ---------------------------------------------------------------------------
# Column width tweaks that accompanied the status text getting longer
# ("In Translation" -> "Ready for handoff" widened the Status/priority-date
# columns). ColumnWidth is quantized internally to the nearest 1/6 character,
# so we pick the input value that rounds to the closest achievable width.
# ---------------------------------------------------------------------------
$overview.Columns.Item(5).ColumnWidth = 16.333333
$overview.Columns.Item(6).ColumnWidth = 16.333333
$zhcn.Columns.Item(3).ColumnWidth = 16.333333
$dede.Columns.Item(3).ColumnWidth = 16.333333
